$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Data edits (applied using the ORIGINAL column layout, before column A is removed) ---

# "Envoi mails" (H2): "Oui" -> boolean TRUE
$ws.Range("H2").Value = $true

# "Montant adhésion année n" (N2): 10 -> 0
$ws.Range("N2").Value = 0

# "Total année n" (Q2): 10 -> 0
$ws.Range("Q2").Value = 0

# "Remarques" (T2): "/" -> cleared
$ws.Range("T2").ClearContents()

# --- Remove the "Enregistrement" column (column A), shifting everything one column left ---
$ws.Columns.Item(1).Delete()

# --- Register the new date/time number format (yyyy-mm-dd h:mm:ss) used for upcoming PDF
#     auto-fill work, without leaving any visible cell using it yet. Apply it to a scratch
#     cell far outside the used range, then remove that row so the style definition is kept
#     but the sheet's dimension / visible data stay untouched. ---
$ws.Range("Z99").NumberFormat = "yyyy-mm-dd h:mm:ss"
$ws.Rows.Item(99).Delete()
